$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$values = @{
    "H11" = 28571662
    "H113" = 100000
    "H137" = 1052.7812
    "H138" = 4010.0557
    "H141" = 2156743.5
    "H15" = 3635.6738
    "H33" = 93.42856999999999
    "H34" = 10000
    "H36" = 10000
    "H64" = 3666.6667
    "H67" = 3666.6667
    "I11" = 28571662
    "I137" = 799.5925999999999
    "I141" = 3501209
    "I15" = 3635.6738
    "I64" = 3000
    "I67" = 3000
    "J113" = 0
    "J138" = 3666.3333
    "J33" = 115.833336
    "J34" = 10000
    "J36" = 10000
    "K11" = 28571662
    "K137" = 2398.7778
    "K141" = 10503627
    "K15" = 10907.0214
    "K64" = 3000
    "K67" = 3000
    "L113" = 0
    "L138" = 10998.9999
    "L33" = 115.833336
    "L34" = 10000
    "L36" = 10000
    "M11" = -28571522
    "M137" = 151.2222000000002
    "M141" = -10498447
    "M15" = -10738.0214
    "M64" = -2752
    "M67" = -2142
    "N138" = -21278.9999
    "N33" = -573.833336
    "N34" = -10406
    "N36" = -11430
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$clearCells = @("N113")
foreach ($ref in $clearCells) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$values = @{
    "H102" = 927.8333
    "H122" = 1564.2
    "H125" = 57142
    "H132" = 1390.2549
    "H136" = 2618.0833
    "H32" = 3019.1455
    "H45" = 1311.1666
    "H61" = 2618.0833
    "H74" = 1074.4595
    "H77" = 1074.4595
    "I102" = 927.8333
    "I122" = 1557.4762
    "I132" = 908.75
    "I136" = 1805.5358
    "I32" = 2601.1042
    "I45" = 1067.15
    "I61" = 1805.5358
    "I74" = 858.1852
    "I77" = 858.1852
    "J122" = 1599.5
    "J125" = 57142
    "J136" = 5462
    "J61" = 5462
    "J74" = 1658.4
    "J77" = 1658.4
    "K102" = 927.8333
    "K122" = 4672.4286
    "K132" = 2726.25
    "K136" = 5416.607400000001
    "K32" = 2601.1042
    "K45" = 1067.15
    "K61" = 1805.5358
    "K74" = 858.1852
    "K77" = 4290.926
    "L122" = 4798.5
    "L125" = 57142
    "L136" = 16386
    "L61" = 5462
    "L74" = 1658.4
    "L77" = 8292
    "M102" = 694.1667
    "M122" = -2222.4286
    "M132" = -196.25
    "M136" = -2866.607400000001
    "M32" = -2314.1042
    "M45" = -690.1500000000001
    "M61" = -1593.5358
    "M74" = 15.81479999999999
    "M77" = 77.07399999999961
    "N122" = -9698.5
    "N125" = -66982
    "N136" = -21486
    "N61" = -5886
    "N74" = -3406.4
    "N77" = -17028
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$values = @{
    "H134" = 4306.1226
    "H20" = 3150.6924
    "H22" = 301
    "H94" = 610.0833
    "I134" = 4873.4595
    "I20" = 2849.889
    "I22" = 301
    "I94" = 652.625
    "K134" = 14620.3785
    "K20" = 2849.889
    "K22" = 301
    "K94" = 652.625
    "M134" = -12085.3785
    "M20" = -2602.889
    "M22" = -128
    "M94" = -201.625
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$values = @{
    "H122" = 3313.077
    "H132" = 1204.8445
    "H31" = 2141.276
    "H34" = 2141.276
    "H62" = 3250
    "H65" = 3250
    "I122" = 1752.75
    "I132" = 798.6389
    "I31" = 2229.077
    "I34" = 2229.077
    "J132" = 2829.6667
    "K122" = 5258.25
    "K132" = 2395.9167
    "K31" = 2229.077
    "K34" = 2229.077
    "L132" = 8489.000100000001
    "M122" = -2808.25
    "M132" = 134.0832999999998
    "M31" = -1934.077
    "M34" = -2027.077
    "N132" = -13549.0001
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$values = @{
    "H122" = 832.7143
    "H131" = 13422.516
    "H2" = 94.23077000000001
    "H3" = 0
    "H36" = 1217.3334
    "H99" = 1690
    "I3" = 0
    "I36" = 1217.3334
    "I99" = 380
    "J122" = 906.125
    "J131" = 14257.184
    "J2" = 40
    "J3" = 0
    "J99" = 3000
    "K3" = 0
    "K36" = 3652.0002
    "K99" = 1140
    "L122" = 8155.125
    "L131" = 42771.552
    "L2" = 240
    "L3" = 0
    "L99" = 9000
    "M36" = -3483.0002
    "M99" = 1106
    "N122" = -13055.125
    "N131" = -52851.552
    "N2" = -466
    "N99" = -13492
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$clearCells = @("M3", "N3")
foreach ($ref in $clearCells) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$values = @{
    "H102" = 3667.077
    "H122" = 1374.2632
    "H132" = 741610.4399999999
    "H80" = 3750
    "H83" = 3750
    "I102" = 3556
    "I122" = 1187.9231
    "I132" = 1013505.4
    "J102" = 5000
    "J122" = 1778
    "K102" = 3556
    "K122" = 3563.7693
    "K132" = 3040516.2
    "L102" = 5000
    "L122" = 5334
    "M102" = -1934
    "M122" = -1113.7693
    "M132" = -3037986.2
    "N102" = -8244
    "N122" = -10234
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$values = @{
    "H122" = 6068.1
    "H132" = 1198.6154
    "H133" = 69663
    "H136" = 2517.3823
    "H22" = 5350
    "H27" = 5350
    "H46" = 2333.3333
    "H55" = 489.18182
    "H82" = 1380.8462
    "H85" = 1380.8462
    "I122" = 5085.125
    "I132" = 825.4792
    "I136" = 1446.4166
    "I22" = 5350
    "I27" = 5350
    "I46" = 2000
    "I55" = 434.875
    "I82" = 1155.2
    "I85" = 1155.2
    "J133" = 69663
    "J22" = 0
    "J27" = 0
    "J46" = 2500
    "J55" = 634
    "K122" = 15255.375
    "K132" = 2476.4376
    "K136" = 4339.2498
    "K22" = 5350
    "K27" = 5350
    "K46" = 2000
    "K55" = 434.875
    "K82" = 1155.2
    "K85" = 1155.2
    "L133" = 69663
    "L22" = 0
    "L27" = 0
    "L46" = 2500
    "L55" = 634
    "M122" = -12805.375
    "M132" = 53.5623999999998
    "M136" = -1789.2498
    "M22" = -5055
    "M27" = -5243
    "M46" = -1812
    "M55" = -261.875
    "M82" = -794.2
    "M85" = 92.79999999999995
    "N133" = -74723
    "N46" = -2876
    "N55" = -980
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$clearCells = @("N22", "N27")
foreach ($ref in $clearCells) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$values = @{
    "H122" = 44244.668
    "I122" = 52895.332
    "K122" = 158685.996
    "M122" = -156235.996
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
